# Update cryptos price (D) and 1h volume-change (E) columns to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.04"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.631.84"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'212.42"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.252"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").Value = "'0.0624"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'19.02"
$ws.Range("E10").Value = "  +2.73%  "
$ws.Range("D11").Value = "'0.0838"
$ws.Range("E11").Value = "  +3.04%  "
$ws.Range("D12").Value = "1.858.46"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "1.633.20"
$ws.Range("E13").Value = "  +1.55%  "
$ws.Range("D14").Value = "'4.08"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "'0.527"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "26.634.39"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "'63.03"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "'209.24"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'9.42"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").Value = "  +3.96%  "
$ws.Range("D25").Value = "'146.92"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'6.81"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "'15.35"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'0.0521"
$ws.Range("E30").Value = "  +5.10%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "'2.95"
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "1.168.11"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "'0.0168"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "'0.504"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +0.93%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "1.768.77"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").Value = "'91.95"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "'54.67"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("D48").Value = "'0.0510"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("E49").Value = "  +4.45%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("E51").Value = "  +0.07%  "
